# Update crypto price/volume data as per the latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.096.68"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "2.757.74"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.31"
$ws.Range("E5").Value = "  -2.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.91"
$ws.Range("E6").Value = "  -1.42%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -3.21%  "
$ws.Range("E9").Value = "  -4.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.89"
$ws.Range("E10").Value = "  -13.33%  "
$ws.Range("E11").Value = "  +3.40%  "
$ws.Range("E12").Value = "  -3.47%  "
$ws.Range("D13").Value = "3.246.96"
$ws.Range("E13").Value = "  -0.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.96"
$ws.Range("E14").Value = "  -2.11%  "
$ws.Range("D15").Value = "63.841.97"
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("E16").Value = "  -5.53%  "
$ws.Range("D17").Value = "2.763.24"
$ws.Range("E17").Value = "  -1.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.13"
$ws.Range("E18").Value = "  -2.40%  "
$ws.Range("E19").Value = "  -5.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "359.14"
$ws.Range("E20").Value = "  -2.31%  "
$ws.Range("E21").Value = "  -6.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.527"
$ws.Range("E23").Value = "  -8.66%  "
$ws.Range("E24").Value = "  -3.80%  "
$ws.Range("E25").Value = "  -3.74%  "
$ws.Range("E26").Value = "  -3.85%  "
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("D28").Value = "0.0₃0903"
$ws.Range("E28").Value = "  -7.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.34"
$ws.Range("E29").Value = "  +0.80%  "
$ws.Range("E30").Value = "  -4.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.35"
$ws.Range("E31").Value = "  +6.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "169.71"
$ws.Range("E32").Value = "  -1.80%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.20"
$ws.Range("E33").Value = "  -3.42%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.93"
$ws.Range("E34").Value = "  -5.95%  "
$ws.Range("E35").Value = "  -2.04%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("E37").Value = "  -2.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  -2.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "347.01"
$ws.Range("E39").Value = "  +1.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.31"
$ws.Range("E40").Value = "  +0.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.18"
$ws.Range("E41").Value = "  -2.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.10"
$ws.Range("E42").Value = "  -3.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.55"
$ws.Range("E43").Value = "  -4.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.93"
$ws.Range("E44").Value = "  -2.87%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0588"
$ws.Range("E45").Value = "  -3.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "137.37"
$ws.Range("E46").Value = "  -0.98%  "
$ws.Range("E47").Value = "  -3.97%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.627"
$ws.Range("E48").Value = "  -4.03%  "
$ws.Range("E49").Value = "  -2.49%  "
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.06"
$ws.Range("E51").Value = "  +0.26%  "
